# Adapt column header formatting to respective input file names (#7)
#   "<header>_old" -> "<header>_FV2410"
#   "<header>_new" -> "<header>_FV2504"
# and turn the data range into a proper Excel Table with the header
# row frozen, as produced by the corresponding xlsx export change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row values: "_old" suffix -> "_FV2410", "_new" suffix -> "_FV2504"
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2. Freeze the header row (top row stays fixed while scrolling)
$win = $excel.ActiveWindow
[void]$ws.Range("A2").Select()
[void]($win.FreezePanes = $true)

# 3. Turn the data range into an Excel Table ("Table1") using the renamed headers
$range = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"

Write-Output "done"
